# Re-pair rows whose data got shuffled between "observation" records.
# Each affected row's entire A:AY content is swapped/rotated with another
# row's content (full-row moves; row numbers/ordering in the sheet stay
# put, only the payload each row carries changes).
#
# Mapping (destination row -> row whose current content it should receive):
#   10<-11  11<-10
#   12<-13  13<-12
#   18<-19  19<-20  20<-21  21<-18
#   22<-23  23<-22
#   27<-28  28<-29  29<-30  30<-27
#   31<-34  32<-31  33<-32  34<-33
#   35<-36  36<-37  37<-38  38<-35
#   46<-47  47<-46
#   51<-54  52<-51  53<-52  54<-53
#   55<-56  56<-55

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceOf = @{
    10 = 11; 11 = 10;
    12 = 13; 13 = 12;
    18 = 19; 19 = 20; 20 = 21; 21 = 18;
    22 = 23; 23 = 22;
    27 = 28; 28 = 29; 29 = 30; 30 = 27;
    31 = 34; 32 = 31; 33 = 32; 34 = 33;
    35 = 36; 36 = 37; 37 = 38; 38 = 35;
    46 = 47; 47 = 46;
    51 = 54; 52 = 51; 53 = 52; 54 = 53;
    55 = 56; 56 = 55
}

# Copy in three column bands (A:X, Z:Z, AB:AY) and deliberately skip Y and
# AA ("Startdatum"/"Slutdatum"): every affected row carries the identical
# "2026-01-18" literal there, so the band is a no-op for this edit, and
# round-tripping a date-shaped string through Range.Value2 makes Excel's
# COM layer re-type it as a real date serial (with a new number format) -
# not the plain inline-string cell the source file has. Skipping the
# unchanged band sidesteps that re-typing entirely.
$bands = @("A:X", "Z:Z", "AB:AY")

# Read every source row's payload for each band BEFORE writing anything
# back, so rows that both give and receive data in the same pass never
# clobber a value that still needs to be read.
$snapshot = @{}
foreach ($destRow in $sourceOf.Keys) {
    $srcRow = $sourceOf[$destRow]
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowData = @{}
        foreach ($band in $bands) {
            $parts = $band.Split(":")
            $addr = "$($parts[0])$($srcRow):$($parts[1])$($srcRow)"
            $rowData[$band] = $ws.Range($addr).Value2
        }
        $snapshot[$srcRow] = $rowData
    }
}

foreach ($destRow in $sourceOf.Keys) {
    $srcRow = $sourceOf[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($band in $bands) {
        $parts = $band.Split(":")
        $addr = "$($parts[0])$($destRow):$($parts[1])$($destRow)"
        $ws.Range($addr).Value2 = $rowData[$band]
    }
}
